# Generate Report for Handback
# Adds a new handback row (row 4) for file
# "e013cfbb-0a7b-4e95-9277-c9ce08160572.md" to the Overview sheet plus the
# per-locale zh-cn and de-de sheets, including hyperlinks and table resizes.

$wb = $excel.ActiveWorkbook

$fileBase   = "e013cfbb-0a7b-4e95-9277-c9ce08160572"
$mdFile     = "$fileBase.md"
$pathName   = "e2e\$fileBase.md"
$ext        = ".md"
$statusSync = "Handed back: in sync with en-US"
$srcPath    = "e2e"
$priority   = "ht"

$xlfZh   = "$fileBase.1c835c3700a86dd172b2a672485b7c1b5a974756.zh-cn.xlf"
$xlfDe   = "$fileBase.1c835c3700a86dd172b2a672485b7c1b5a974756.de-de.xlf"

$hoDateZh  = "2016-09-04 16:49:07"
$hbDateZh  = "2016-09-04 16:49:31"
$hoDateDe  = "2016-09-04 16:49:12"
$hbDateDe  = "2016-09-04 16:49:39"
$latestDate = "2016-09-04 16:49:12"

# ----------------------------------------------------------------------
# Sheet "Overview" -> row 4
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdFile
$wsOverview.Range("B4").Value = $pathName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e013cfbb0a7b4e959277c9ce08160572handback/e2e/$mdFile", "", "", $pathName)
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").Value = $latestDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ----------------------------------------------------------------------
# Sheet "zh-cn" -> row 4
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdFile
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e013cfbb0a7b4e959277c9ce08160572handback/e2e/$mdFile", "", "", $mdFile)
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $statusSync
$wsZh.Range("D4").Value = $srcPath
$wsZh.Range("E4").Value = $priority
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $xlfZh
$wsZh.Range("H4").Value = $hoDateZh
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $mdFile
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e013cfbb0a7b4e959277c9ce08160572handback/e2e/$mdFile", "", "", $mdFile)
$wsZh.Range("J4").Value = $xlfZh
$wsZh.Range("K4").Value = $hbDateZh
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ----------------------------------------------------------------------
# Sheet "de-de" -> row 4
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdFile
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e013cfbb0a7b4e959277c9ce08160572handback/e2e/$mdFile", "", "", $mdFile)
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $statusSync
$wsDe.Range("D4").Value = $srcPath
$wsDe.Range("E4").Value = $priority
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $xlfDe
$wsDe.Range("H4").Value = $hoDateDe
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $mdFile
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e013cfbb0a7b4e959277c9ce08160572handback/e2e/$mdFile", "", "", $mdFile)
$wsDe.Range("J4").Value = $xlfDe
$wsDe.Range("K4").Value = $hbDateDe
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))

Write-Host "Handback row added for $mdFile"
